$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(17).Insert()
$ws.Range("B16:J16").Copy($ws.Range("B17:J17"))

Write-Host "done"
